$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The three data rows (2, 3, 4) are cyclically shifted:
#   new row 2 <- old row 3's data
#   new row 3 <- old row 4's data
#   new row 4 <- old row 2's data
# Only the "Id" (A), "Ost" (Q) and "Nord" (R) values actually differ between the
# rows, plus a handful of always-blank marker cells (J, K, L, N, AF) that are
# present on some rows and absent on others. Apply the value swaps directly and
# add/remove the blank marker cells to match.

# --- A (Id) column ---
$ws.Range("A2").Value = 104789708
$ws.Range("A3").Value = 104789709
$ws.Range("A4").Value = 104793729

# --- Q (Ost) column ---
$ws.Range("Q2").Value = 754058.7871261307
$ws.Range("Q3").Value = 754080.4528164999
$ws.Range("Q4").Value = 754004.464596002

# --- R (Nord) column ---
$ws.Range("R2").Value = 7197185.128098626
$ws.Range("R3").Value = 7197162.135852392
$ws.Range("R4").Value = 7197229.324461597

# --- Blank marker cells (J2, K2, L2, N2, AF2) now appear on row 2 ---
# (row 3 already has them and keeps them unchanged; row 4 loses them below)
$ws.Range("J2").Font.Bold = $false
$ws.Range("K2").Font.Bold = $false
$ws.Range("L2").Font.Bold = $false
$ws.Range("N2").Font.Bold = $false
$ws.Range("AF2").Font.Bold = $false

# --- Blank marker cells (J4, K4, L4, N4, AF4) are removed from row 4 ---
$ws.Range("J4").ClearContents()
$ws.Range("K4").ClearContents()
$ws.Range("L4").ClearContents()
$ws.Range("N4").ClearContents()
$ws.Range("AF4").ClearContents()
